# Add the new "Canada" source-IP row (row 35) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(35, 1).Value = "70.26.209.190"
$ws.Cells.Item(35, 2).Value = "Vaureuil, Canada"

# Match the author's final selection/scroll state as closely as the
# object model allows.
$excel.ActiveWindow.ScrollRow = 29
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A35").Select()
